$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "~0" shared-string placeholders in the p-value columns are replaced
# with plain numeric 0s.
$ws.Range("C5").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0

# Add a bold footnote below the table explaining what "0" means, reusing
# the now-unused shared string that used to hold "~0".
$ws.Range("B18").Value = "0 refers to tending to zero"
$ws.Range("B18").Font.Bold = $true

# Apply a green-yellow-red 3-color scale over the p-value columns.
$rng = $ws.Range("C4:D16")
$cf = $rng.FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria(1).FormatColor.Color = 8109667
$cf.ColorScaleCriteria(3).FormatColor.Color = 7039480

$ws.Range("B18").Select()
